$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-05-12"

# Update the row label for May in column A (row 6)
$ws.Range("A6").Value = "May (through 05-12)"

# Update May row (row 6) values for years 2015-2022 (columns B-I)
$ws.Range("B6").Value = 9
$ws.Range("C6").Value = 15
$ws.Range("D6").Value = 25
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 13
$ws.Range("G6").Value = 24
$ws.Range("H6").Value = 45
$ws.Range("I6").Value = 42

# Update Total row (row 7) values for years 2015-2022 (columns B-I)
$ws.Range("B7").Value = 98
$ws.Range("C7").Value = 177
$ws.Range("D7").Value = 278
$ws.Range("E7").Value = 261
$ws.Range("F7").Value = 168
$ws.Range("G7").Value = 286
$ws.Range("H7").Value = 568
$ws.Range("I7").Value = 594
